# Fix contact information missing from short resumes.
#
# Insert a new centered paragraph containing the contact information
# directly after the "Dheeraj Chand" name paragraph and before the
# "PROFESSIONAL SUMMARY" heading.
#
# We do this via Find/Replace (embedding a "^p" paragraph-mark in the
# replacement text) rather than Paragraphs/Range InsertParagraphAfter,
# because InsertParagraphAfter/Before on an adjacent paragraph causes the
# new paragraph to inherit unwanted formatting (bold/28pt run formatting
# from the name line, or the Heading2 paragraph style from the summary
# heading). Performing the split inside the "Dheeraj Chand" run via
# Find/Replace keeps the new paragraph's run free of any rPr and its
# pPr limited to just the centered alignment, matching a plain,
# unformatted contact line.

$d = $word.ActiveDocument

$null = $d.Content.Find.Execute(
    "Dheeraj Chand",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)

# Ensure the newly created contact-info paragraph (now paragraph 2) is
# centered, consistent with the name paragraph above it.
$contactPara = $d.Paragraphs.Item(2)
$contactPara.Alignment = 1
